# Update column F (dSF) values to match repulled/pushed data and mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 0
$ws.Range("F8").Value = -3
$ws.Range("F10").Value = 9
$ws.Range("F11").Value = -9
$ws.Range("F12").Value = -5
$ws.Range("F18").Value = -3
$ws.Range("F21").Value = -1
$ws.Range("F22").Value = 6
